$wb = $excel.ActiveWorkbook

# Insert a new worksheet named "Descriptors" right after the "QuantityValue" sheet
# and before the "Acquisition" sheet.
$afterSheet = $wb.Worksheets.Item("QuantityValue")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "Descriptors"

$newSheet.Range("A1").Value = "descriptor_name"
$newSheet.Range("B1").Value = "descriptor_thing"
